$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text formatting (values contain
# "." as thousands separators and must not be auto-converted to numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '63.270.95'
$ws.Range('E2').Value = '  -3.35%  '
$ws.Range('D3').Value = '3.257.03'
$ws.Range('E3').Value = '  -4.20%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = '175.97'
$ws.Range('E5').Value = '  -3.26%  '
$ws.Range('D6').Value = '521.30'
$ws.Range('E6').Value = '  -1.68%  '
$ws.Range('D7').Value = '0.593'
$ws.Range('E7').Value = '  -2.82%  '
$ws.Range('D8').Value = '3.254.38'
$ws.Range('E8').Value = '  -3.97%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  -3.54%  '
$ws.Range('D11').Value = '53.07'
$ws.Range('E11').Value = '  -7.48%  '
$ws.Range('E12').Value = '  -2.23%  '
$ws.Range('E13').Value = '  -0.96%  '
$ws.Range('E14').Value = '  -3.93%  '
$ws.Range('D15').Value = '3.777.28'
$ws.Range('E15').Value = '  -4.30%  '
$ws.Range('D16').Value = '3.256.99'
$ws.Range('E16').Value = '  -4.66%  '
$ws.Range('D17').Value = '0.116'
$ws.Range('E17').Value = '  -6.17%  '
$ws.Range('D18').Value = '63.246.83'
$ws.Range('E18').Value = '  -3.18%  '
$ws.Range('D19').Value = '17.22'
$ws.Range('E19').Value = '  -1.82%  '
$ws.Range('D20').Value = '11.03'
$ws.Range('E20').Value = '  -1.82%  '
$ws.Range('D21').Value = '0.963'
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('D22').Value = '367.05'
$ws.Range('E22').Value = '  -2.60%  '
$ws.Range('D23').Value = '3.74'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').Value = '80.54'
$ws.Range('E24').Value = '  -2.81%  '
$ws.Range('D25').Value = '10.99'
$ws.Range('E25').Value = '  +1.53%  '
$ws.Range('D26').Value = '3.91'
$ws.Range('E26').Value = '  +9.29%  '
$ws.Range('E27').Value = '  +4.69%  '
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('D29').Value = '11.26'
$ws.Range('E29').Value = '  -2.97%  '
$ws.Range('D30').Value = '8.19'
$ws.Range('E30').Value = '  -3.84%  '
$ws.Range('D31').Value = '656.58'
$ws.Range('E31').Value = '  -3.42%  '
$ws.Range('D32').Value = '28.37'
$ws.Range('E32').Value = '  -4.76%  '
$ws.Range('D33').Value = '6.38'
$ws.Range('E33').Value = '  -4.59%  '
$ws.Range('D34').Value = '11.14'
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('D35').Value = '0.104'
$ws.Range('E35').Value = '  -0.97%  '
$ws.Range('D36').Value = '57.57'
$ws.Range('E36').Value = '  -6.43%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('D38').Value = '36.54'
$ws.Range('E38').Value = '  -0.31%  '
$ws.Range('D39').Value = '0.375'
$ws.Range('E39').Value = '  -2.61%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0714'
$ws.Range('E41').Value = '  +14.06%  '
$ws.Range('E42').Value = '  -3.43%  '
$ws.Range('D43').Value = '2.889.11'
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('D44').Value = '2.48'
$ws.Range('E44').Value = '  +6.55%  '
$ws.Range('E45').Value = '  -1.34%  '
$ws.Range('D46').Value = '0.0391'
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = '2.81'
$ws.Range('E47').Value = '  +9.83%  '
$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D48').Value = '2.57'
$ws.Range('E48').Value = '  -7.15%  '
$ws.Range('D49').Value = '2.99'
$ws.Range('E49').Value = '  +5.63%  '
$ws.Range('D50').Value = '135.72'
$ws.Range('E50').Value = '  +0.34%  '
$ws.Range('E51').Value = '  -2.11%  '
